$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Austria (row 20) case figures ---
$ws.Range("B20").Value = 14499
$ws.Range("C20").Value = 23
$ws.Range("E20").Value = 5103

# --- Update Kazajistan (row 67) case figures ---
$ws.Range("B67").Value = 1480
$ws.Range("C67").Value = 78
$ws.Range("E67").Value = 1178

# --- Insert "Consejo Danes para los Refugiados" between Montenegro and Isla de Man ---
# This shifts Isla de Man / Vietnam down one row (Sri Lanka row stays put because the
# old "Consejo Danes para los Refugiados" row that used to sit just above it is gone).

# Row 114 becomes the new "Consejo Danes para los Refugiados" entry
$ws.Range("A114").Value = "Consejo Danes para los Refugiados"
$ws.Range("B114").Value = 287
$ws.Range("C114").Value = 20
$ws.Range("D114").Value = 25
$ws.Range("E114").Value = 239
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 1
$ws.Range("H114").Value = 23

# Row 115 becomes "Isla de Man" (previously row 114's data)
$ws.Range("A115").Value = "Isla de Man"
$ws.Range("B115").Value = 284
$ws.Range("C115").Value = 0
$ws.Range("D115").Value = 154
$ws.Range("E115").Value = 126
$ws.Range("F115").Value = 13
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 4

# Row 116 becomes "Vietnam" (previously row 115's data)
$ws.Range("A116").Value = "Vietnam"
$ws.Range("B116").Value = 268
$ws.Range("C116").Value = 0
$ws.Range("D116").Value = 194
$ws.Range("E116").Value = 74
$ws.Range("F116").Value = 8
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 0

# Row 117 ("Sri Lanka") is unchanged - already correct.
